$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (46081 -> 46082) for every data row (rows 2 through 95).
$lastRow = 95
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46081) {
        $cell.Value = 46082
    }
}
